$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 42 (the "004870019" row) for the THEOMAR
# account, which now has a higher Saldo (542.71) and must move up the
# list to keep it sorted descending by Saldo.
$ws.Rows.Item(42).Insert()

$ws.Cells.Item(42, 1).Value = "'004382374"
$ws.Cells.Item(42, 2).Value = "THEOMAR"
$ws.Cells.Item(42, 3).Value = 542.71

# The original THEOMAR row (previously row 382, now shifted down to 383
# because of the insert above) is removed from its old sorted position.
$ws.Rows.Item(383).Delete()
